$d = $word.ActiveDocument

# Find the range spanning the 7 "Author"-styled paragraphs (Margaret Westbury
# through Curtis Sharma) and delete them entirely, including their paragraph
# marks, so the Subtitle paragraph is followed directly by the Date paragraph.

$startPara = $null
$endPara = $null

foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Author") {
        if ($startPara -eq $null) {
            $startPara = $p
        }
        $endPara = $p
    }
}

if ($startPara -ne $null -and $endPara -ne $null) {
    $rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $rng.Delete()
}
